$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new log entry row 41 ---
# Copy existing formatting from the row above (row 40) as a starting point,
# then set the new values and tweak the look of the new row.

$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B40").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("C40").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A41").Value = 45742
$ws.Range("B41").Value = 4
$ws.Range("C41").Value = "Worked on transforming the final output format"

# Center the hours cell and mark its fill as explicitly set
$ws.Range("B41").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B41").Interior.ColorIndex = -4142   # xlColorIndexNone

# Description cell: general (left) horizontal alignment, vertical centered,
# wrap text, only a thin top/bottom border (no left/right), fill explicitly
# marked as none as well.
$ws.Range("C41").HorizontalAlignment = 1       # xlGeneral
$ws.Range("C41").VerticalAlignment = -4108     # xlCenter
$ws.Range("C41").WrapText = $true
$ws.Range("C41").Interior.ColorIndex = -4142   # xlColorIndexNone
$ws.Range("C41").Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> none
$ws.Range("C41").Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none
$ws.Range("C41").Borders.Item(8).LineStyle = 1       # xlEdgeTop -> thin
$ws.Range("C41").Borders.Item(8).Weight = 2
$ws.Range("C41").Borders.Item(9).LineStyle = 1       # xlEdgeBottom -> thin
$ws.Range("C41").Borders.Item(9).Weight = 2

# --- Update the selection to mirror where the user left the cursor ---
$ws.Range("C44").Select()

Write-Host "done"
